# Apply crypto price / volume(1h) updates to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'61.621.62"
$ws.Range('D3').Value = "'2.895.92"
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'568.76"
$ws.Range('E5').Value = '  -4.37%  '
$ws.Range('D6').Value = "'144.40"
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('D9').Value = "'2.895.06"
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('D10').Value = "'7.02"
$ws.Range('E10').Value = '  -3.61%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').Value = "'32.05"
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = "'3.375.72"
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').Value = "'61.603.40"
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').Value = "'2.891.73"
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = "'433.04"
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('D21').Value = "'13.12"
$ws.Range('E21').Value = '  -2.69%  '
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('D24').Value = "'79.37"
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('D26').Value = "'10.02"
$ws.Range('E26').Value = '  -10.89%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('D30').Value = "'7.04"
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  -6.54%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = "'0.106"
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').Value = "'25.52"
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('D36').Value = "'0.961"
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('E37').Value = '  -3.32%  '
$ws.Range('D38').Value = "'48.89"
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('E39').Value = '  -5.38%  '
$ws.Range('E40').Value = '  -10.09%  '
$ws.Range('E41').Value = '  -2.04%  '
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('D43').Value = "'39.65"
$ws.Range('E43').Value = '  +1.88%  '
$ws.Range('D44').Value = "'0.267"
$ws.Range('E44').Value = '  -5.07%  '
$ws.Range('D45').Value = "'2.708.99"
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('D46').Value = "'133.00"
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').Value = "'0.0335"
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = "'346.80"
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = "'21.63"
$ws.Range('E51').Value = '  -5.23%  '
